$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1984536082474227
$ws.Range("C2").Value = 0.5309278350515464
$ws.Range("J2").Value = 0.03350515463917526
$ws.Range("P2").Value = 0.1494845360824742
$ws.Range("S2").Value = 0.08762886597938144
$ws.Range("B3").Value = 0.004672897196261682
$ws.Range("C3").Value = 0.03271028037383177
$ws.Range("J3").Value = 0.04205607476635514
$ws.Range("P3").Value = 0.7102803738317757
$ws.Range("S3").Value = 0.2102803738317757
$ws.Range("J4").Value = 0.06818181818181818
$ws.Range("O4").Value = 0.02272727272727273
$ws.Range("P4").Value = 0.6363636363636364
$ws.Range("S4").Value = 0.2727272727272727
$ws.Range("B6").Value = 0.07174887892376682
$ws.Range("D6").Value = 0.0179372197309417
$ws.Range("F6").Value = 0.05381165919282511
$ws.Range("J6").Value = 0.2780269058295964
$ws.Range("O6").Value = 0.0179372197309417
$ws.Range("Q6").Value = 0.1704035874439462
$ws.Range("R6").Value = 0.05829596412556054
$ws.Range("S6").Value = 0.3318385650224215
$ws.Range("B7").Value = 0.1611111111111111
$ws.Range("D7").Value = 0.005555555555555556
$ws.Range("E7").Value = 0.01111111111111111
$ws.Range("F7").Value = 0.02777777777777778
$ws.Range("J7").Value = 0.1666666666666667
$ws.Range("O7").Value = 0.005555555555555556
$ws.Range("Q7").Value = 0.1388888888888889
$ws.Range("R7").Value = 0.09444444444444444
$ws.Range("S7").Value = 0.3888888888888889
$ws.Range("B8").Value = 0.116710875331565
$ws.Range("D8").Value = 0.01856763925729443
$ws.Range("F8").Value = 0.07161803713527852
$ws.Range("J8").Value = 0.1352785145888594
$ws.Range("O8").Value = 0.007957559681697613
$ws.Range("Q8").Value = 0.1962864721485411
$ws.Range("R8").Value = 0.09549071618037135
$ws.Range("S8").Value = 0.3580901856763926
$ws.Range("B9").Value = 0.1451612903225807
$ws.Range("D9").Value = 0.01612903225806452
$ws.Range("F9").Value = 0.06854838709677419
$ws.Range("J9").Value = 0.08870967741935484
$ws.Range("O9").Value = 0.01612903225806452
$ws.Range("Q9").Value = 0.157258064516129
$ws.Range("R9").Value = 0.1169354838709677
$ws.Range("S9").Value = 0.3911290322580645
$ws.Range("B10").Value = 0.1384959046909903
$ws.Range("D10").Value = 0.02159344750558451
$ws.Range("E10").Value = 0.0007446016381236039
$ws.Range("F10").Value = 0.05956813104988831
$ws.Range("J10").Value = 0.1422189128816083
$ws.Range("O10").Value = 0.01116902457185406
$ws.Range("Q10").Value = 0.197319434102755
$ws.Range("R10").Value = 0.08116157855547282
$ws.Range("S10").Value = 0.347728965003723
$ws.Range("F11").Value = 0.003787878787878788
$ws.Range("G11").Value = 0.1590909090909091
$ws.Range("J11").Value = 0.06818181818181818
$ws.Range("K11").Value = 0.1893939393939394
$ws.Range("L11").Value = 0.5454545454545454
$ws.Range("S11").Value = 0.03409090909090909
$ws.Range("G12").Value = 0.7707006369426752
$ws.Range("J12").Value = 0.1847133757961783
$ws.Range("K12").Value = 0.006369426751592357
$ws.Range("L12").Value = 0.01273885350318471
$ws.Range("S12").Value = 0.02547770700636943
$ws.Range("G13").Value = 0.6097560975609756
$ws.Range("J13").Value = 0.3170731707317073
$ws.Range("S13").Value = 0.07317073170731707
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.02777777777777778
$ws.Range("H15").Value = 0.15
$ws.Range("I15").Value = 0.1
$ws.Range("J15").Value = 0.4388888888888889
$ws.Range("K15").Value = 0.04444444444444445
$ws.Range("M15").Value = 0.01666666666666667
$ws.Range("O15").Value = 0.04444444444444445
$ws.Range("S15").Value = 0.1777777777777778
$ws.Range("F16").Value = 0.03418803418803419
$ws.Range("H16").Value = 0.1837606837606838
$ws.Range("I16").Value = 0.09829059829059829
$ws.Range("J16").Value = 0.3632478632478632
$ws.Range("K16").Value = 0.1153846153846154
$ws.Range("M16").Value = 0.008547008547008548
$ws.Range("O16").Value = 0.05128205128205128
$ws.Range("S16").Value = 0.1452991452991453
$ws.Range("F17").Value = 0.009111617312072893
$ws.Range("H17").Value = 0.1526195899772209
$ws.Range("I17").Value = 0.1138952164009112
$ws.Range("J17").Value = 0.469248291571754
$ws.Range("K17").Value = 0.07517084282460136
$ws.Range("M17").Value = 0.01822323462414579
$ws.Range("O17").Value = 0.03644646924829157
$ws.Range("S17").Value = 0.1252847380410023
$ws.Range("F18").Value = 0.01005025125628141
$ws.Range("H18").Value = 0.1155778894472362
$ws.Range("I18").Value = 0.1256281407035176
$ws.Range("J18").Value = 0.4974874371859296
$ws.Range("K18").Value = 0.05527638190954774
$ws.Range("M18").Value = 0.02512562814070352
$ws.Range("N18").Value = 0.005025125628140704
$ws.Range("O18").Value = 0.05527638190954774
$ws.Range("S18").Value = 0.1105527638190955
$ws.Range("F19").Value = 0.02722772277227723
$ws.Range("H19").Value = 0.1765676567656766
$ws.Range("I19").Value = 0.1047854785478548
$ws.Range("J19").Value = 0.3754125412541254
$ws.Range("K19").Value = 0.1023102310231023
$ws.Range("M19").Value = 0.01897689768976898
$ws.Range("O19").Value = 0.07178217821782178
$ws.Range("S19").Value = 0.1229372937293729
